$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Kaspersky Small Office Security price text (shared string used by B19)
$ws.Range("B19").Value = "Kaspersky Small Office Security/$ R$ 1.267"

# C18: update the Kaspersky monthly-equivalent numeric value
$ws.Range("C18").Value = 1267.2

# C19: clear the stale total value (cell becomes blank, keeps its style)
$ws.Range("C19").ClearContents()

# B26: update the budget reserve value
$ws.Range("B26").Value = 5000

# Move the active selection to B26, matching the saved selection state
$ws.Range("B26").Select()
